# Updated symbol list on Sat Dec 17 15:55:59 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates -- force text storage ("@") so that numeric-looking
# strings (e.g. with trailing/leading zeros) are preserved exactly as text,
# matching the workbook's existing inline-string / text cell representation.
$priceUpdates = @{
    "D2"  = "236.94"
    "D3"  = "21.84"
    "D4"  = "5.387"
    "D6"  = "6.474"
    "D7"  = "3.351"
    "D8"  = "0.7991"
    "D9"  = "1.036"
    "D11" = "0.07277"
    "D12" = "0.03125"
    "D13" = "0.02971"
    "D14" = "0.09242"
    "D15" = "0.001659"
    "D16" = "3.264"
    "D17" = "0.04777"
    "D18" = "0.0005713"
    "D19" = "0.006230"
    "D20" = "0.005075"
    "D21" = "0.001050"
    "D22" = "0.0001501"
    "D24" = "3.951"
    "D27" = "0.1061"
    "D40" = "0.04088"
    "D41" = "0.006973"
    "D42" = "0.003502"
    "D43" = "0.1038"
    "D44" = "0.008840"
    "D45" = "0.00005441"
    "D48" = "0.03702"
    "D49" = "0.00002101"
    "D50" = "0.01011"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Volume(1h) column (E) label updates -- plain text, no numeric coercion risk.
$volumeUpdates = @{
    "E18" = "17OneONE"
    "E24" = "23LEOLEO"
    "E42" = "41CEJICEJIBestin24h"
    "E48" = "47BOLOBOLOWorstin24h"
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
